# Auto-generated Excel COM-interop script to apply profit recalculations
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets of the Cerberus_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2784
$ws.Range("J40").Value = 3079.75
$ws.Range("L40").Value = 3079.75
$ws.Range("N40").Value = -3429.75

$ws.Range("H64").Value = 5639.4
$ws.Range("I64").Value = 5639.4
$ws.Range("K64").Value = 5639.4
$ws.Range("M64").Value = -5391.4

$ws.Range("H67").Value = 5639.4
$ws.Range("I67").Value = 5639.4
$ws.Range("K67").Value = 5639.4
$ws.Range("M67").Value = -4781.4

$ws.Range("H103").Value = 484.7
$ws.Range("J103").Value = 711.6
$ws.Range("L103").Value = 2134.8
$ws.Range("N103").Value = -3306.8

$ws.Range("H129").Value = 1733.4286
$ws.Range("I129").Value = 843.0909
$ws.Range("J129").Value = 4998
$ws.Range("K129").Value = 2529.2727
$ws.Range("L129").Value = 14994
$ws.Range("M129").Value = 2470.7273
$ws.Range("N129").Value = -24994

$ws.Range("H132").Value = 4268.4316
$ws.Range("I132").Value = 4080.756
$ws.Range("K132").Value = 12242.268
$ws.Range("M132").Value = -9712.268

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H61").Value = 12094.0625
$ws.Range("I61").Value = 8762.125
$ws.Range("K61").Value = 8762.125
$ws.Range("M61").Value = -8550.125

$ws.Range("H74").Value = 4598.2915
$ws.Range("I74").Value = 2620
$ws.Range("J74").Value = 6011.357
$ws.Range("K74").Value = 2620
$ws.Range("L74").Value = 6011.357
$ws.Range("M74").Value = -1746
$ws.Range("N74").Value = -7759.357

$ws.Range("H77").Value = 4598.2915
$ws.Range("I77").Value = 2620
$ws.Range("J77").Value = 6011.357
$ws.Range("K77").Value = 13100
$ws.Range("L77").Value = 30056.785
$ws.Range("M77").Value = -8732
$ws.Range("N77").Value = -38792.785

$ws.Range("H132").Value = 1967
$ws.Range("I132").Value = 1542.7778
$ws.Range("K132").Value = 4628.3334
$ws.Range("M132").Value = -2098.3334

$ws.Range("H136").Value = 12094.0625
$ws.Range("I136").Value = 8762.125
$ws.Range("K136").Value = 26286.375
$ws.Range("M136").Value = -23736.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4619.3335
$ws.Range("I99").Value = 4571.75
$ws.Range("K99").Value = 4571.75
$ws.Range("M99").Value = -3073.75

$ws.Range("H134").Value = 7844.6665
$ws.Range("I134").Value = 6640.857
$ws.Range("J134").Value = 10653.556
$ws.Range("K134").Value = 19922.571
$ws.Range("L134").Value = 31960.668
$ws.Range("M134").Value = -17387.571
$ws.Range("N134").Value = -37030.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1425.1428
$ws.Range("I16").Value = 1212.6
$ws.Range("J16").Value = 1956.5
$ws.Range("K16").Value = 1212.6
$ws.Range("L16").Value = 1956.5
$ws.Range("M16").Value = -925.5999999999999
$ws.Range("N16").Value = -2530.5

$ws.Range("H22").Value = 519.2
$ws.Range("I22").Value = 138.5
$ws.Range("J22").Value = 773
$ws.Range("K22").Value = 138.5
$ws.Range("L22").Value = 773
$ws.Range("M22").Value = 211.5
$ws.Range("N22").Value = -1473

$ws.Range("H60").Value = 15399
$ws.Range("I60").Value = 5539.7144
$ws.Range("K60").Value = 5539.7144
$ws.Range("M60").Value = -5028.7144

$ws.Range("H113").Value = 1425.1428
$ws.Range("I113").Value = 1212.6
$ws.Range("J113").Value = 1956.5
$ws.Range("K113").Value = 1212.6
$ws.Range("L113").Value = 1956.5
$ws.Range("M113").Value = 957.4000000000001
$ws.Range("N113").Value = -6296.5

$ws.Range("H122").Value = 3908.75
$ws.Range("I122").Value = 3689.4
$ws.Range("J122").Value = 5005.5
$ws.Range("K122").Value = 11068.2
$ws.Range("L122").Value = 15016.5
$ws.Range("M122").Value = -8618.2
$ws.Range("N122").Value = -19916.5

$ws.Range("H132").Value = 2991.6956
$ws.Range("I132").Value = 3222.2856
$ws.Range("K132").Value = 9666.856800000001
$ws.Range("M132").Value = -7136.856800000001

$ws.Range("H134").Value = 4788.9023
$ws.Range("I134").Value = 4283.5293
$ws.Range("K134").Value = 12850.5879
$ws.Range("M134").Value = -10315.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1322.25
$ws.Range("I2").Value = 2540.75
$ws.Range("K2").Value = 15244.5
$ws.Range("M2").Value = -15131.5

$ws.Range("H4").Value = 1220333.5
$ws.Range("I4").Value = 2615575.5
$ws.Range("K4").Value = 7846726.5
$ws.Range("M4").Value = -7846614.5

$ws.Range("H15").Value = 1833.3334
$ws.Range("J15").Value = 2625
$ws.Range("L15").Value = 7875
$ws.Range("N15").Value = -8155

$ws.Range("H34").Value = 2315.2
$ws.Range("I34").Value = 174
$ws.Range("J34").Value = 2850.5
$ws.Range("K34").Value = 522
$ws.Range("L34").Value = 8551.5
$ws.Range("M34").Value = -438
$ws.Range("N34").Value = -8719.5

$ws.Range("H131").Value = 15875144
$ws.Range("I131").Value = 11112619
$ws.Range("J131").Value = 18520992
$ws.Range("K131").Value = 33337857
$ws.Range("L131").Value = 55562976
$ws.Range("M131").Value = -33332817
$ws.Range("N131").Value = -55573056

$ws.Range("H139").Value = 6674783.5
$ws.Range("J139").Value = 12306.533
$ws.Range("L139").Value = 36919.599
$ws.Range("N139").Value = -47199.599

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2278.1765
$ws.Range("I7").Value = 2211.25
$ws.Range("J7").Value = 2438.8
$ws.Range("K7").Value = 2211.25
$ws.Range("L7").Value = 2438.8
$ws.Range("M7").Value = -2099.25
$ws.Range("N7").Value = -2662.8

$ws.Range("H55").Value = 659.0769
$ws.Range("I55").Value = 479.1
$ws.Range("J55").Value = 1259
$ws.Range("K55").Value = 479.1
$ws.Range("L55").Value = 1259
$ws.Range("M55").Value = -306.1
$ws.Range("N55").Value = -1605

$ws.Range("H68").Value = 1927.8462
$ws.Range("I68").Value = 1841.8636
$ws.Range("J68").Value = 2400.75
$ws.Range("K68").Value = 1841.8636
$ws.Range("L68").Value = 2400.75
$ws.Range("M68").Value = -1092.8636
$ws.Range("N68").Value = -3898.75

$ws.Range("H71").Value = 1927.8462
$ws.Range("I71").Value = 1841.8636
$ws.Range("J71").Value = 2400.75
$ws.Range("K71").Value = 9209.318
$ws.Range("L71").Value = 12003.75
$ws.Range("M71").Value = -5465.317999999999
$ws.Range("N71").Value = -19491.75

$ws.Range("H95").Value = 25000.5
$ws.Range("J95").Value = 25000.5
$ws.Range("L95").Value = 25000.5
$ws.Range("N95").Value = -30492.5

$ws.Range("H126").Value = 2278.1765
$ws.Range("I126").Value = 2211.25
$ws.Range("J126").Value = 2438.8
$ws.Range("K126").Value = 6633.75
$ws.Range("L126").Value = 7316.400000000001
$ws.Range("M126").Value = -4163.75
$ws.Range("N126").Value = -12256.4

$ws.Range("H136").Value = 3145.1667
$ws.Range("I136").Value = 1283.3334
$ws.Range("J136").Value = 3610.625
$ws.Range("K136").Value = 3850.0002
$ws.Range("L136").Value = 10831.875
$ws.Range("M136").Value = -1300.0002
$ws.Range("N136").Value = -15931.875

$ws.Range("H140").Value = 53999.2
$ws.Range("J140").Value = 56110.223
$ws.Range("L140").Value = 56110.223
$ws.Range("N140").Value = -66470.223

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4922.778
$ws.Range("I122").Value = 5084.3335
$ws.Range("K122").Value = 15253.0005
$ws.Range("M122").Value = -12803.0005

$ws.Range("H132").Value = 2749.5386
$ws.Range("I132").Value = 2366.875
$ws.Range("K132").Value = 7100.625
$ws.Range("M132").Value = -4570.625
